$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 46: same date-style as prior rows, additional effort of 1h,
# and a comment describing minor documentation/setup changes.
$ws.Range("A45").Copy()
$ws.Range("A46").PasteSpecial(-4122)
$ws.Range("A46").Value = 41234

$ws.Range("B46").Value = 1
$ws.Range("D46").Value = "Minor changes on documentation and setup"

# Match the selection state from the diff (active cell moves to D46)
$ws.Range("D46").Select()
